$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.9
$ws.Range("I2").Value = 1.94
$ws.Range("J2").Value = 3.4
$ws.Range("L2").Value = 1.34
$ws.Range("N2").Value = 3.25
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 1.77
$ws.Range("S2").Value = 3.35
$ws.Range("U2").Value = 1.89
$ws.Range("V2").Value = 2.06
$ws.Range("W2").Value = 1.2
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 21
$ws.Range("AD2").Value = 11
$ws.Range("AG2").Value = 22
$ws.Range("L3").Value = 1.41
$ws.Range("M3").Value = 1.07
$ws.Range("P3").Value = 1.82
$ws.Range("Q3").Value = 1.91
$ws.Range("S3").Value = 3.2
$ws.Range("W3").Value = 1.94
$ws.Range("F4").Value = 1.26
$ws.Range("O4").Value = 1.12
$ws.Range("S4").Value = 2
$ws.Range("F6").Value = 1.81
$ws.Range("G6").Value = 1.93
$ws.Range("H6").Value = 4.4
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 3.75
$ws.Range("L6").Value = 1.38
$ws.Range("N6").Value = 3.75
$ws.Range("R6").Value = 1.37
$ws.Range("T6").Value = 1.8
$ws.Range("V6").Value = 1.25
$ws.Range("Y6").Value = 20
$ws.Range("Z6").Value = 42
$ws.Range("AD6").Value = 22
$ws.Range("AE6").Value = 70
$ws.Range("AF6").Value = 12
$ws.Range("AH6").Value = 21
$ws.Range("AI6").Value = 75
$ws.Range("AJ6").Value = 21
$ws.Range("AO6").Value = 80
$ws.Range("J7").Value = 3.6
$ws.Range("K7").Value = 4.2
$ws.Range("V7").Value = 1.28
$ws.Range("F8").Value = 1.16
$ws.Range("G8").Value = 1.2
$ws.Range("I8").Value = 40
$ws.Range("J8").Value = 7.4
$ws.Range("K8").Value = 10
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 4.5
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 2.26
$ws.Range("Q8").Value = 1.66
$ws.Range("R8").Value = 1.49
$ws.Range("S8").Value = 2.68
$ws.Range("T8").Value = 2.82
$ws.Range("U8").Value = 1.45
$ws.Range("W8").Value = 6
$ws.Range("Y8").Value = 85
$ws.Range("F9").Value = 1.75
$ws.Range("H9").Value = 4.6
$ws.Range("P9").Value = 2.3
$ws.Range("Q9").Value = 1.64
$ws.Range("R9").Value = 1.53
$ws.Range("S9").Value = 2.54
$ws.Range("T9").Value = 1.66
$ws.Range("U9").Value = 2.26
$ws.Range("V9").Value = 1.24
$ws.Range("AF9").Value = 13.5
$ws.Range("G10").Value = 2.72
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.85
$ws.Range("J10").Value = 2.9
$ws.Range("V10").Value = 1.36
$ws.Range("W10").Value = 1.6
$ws.Range("F11").Value = 2.14
$ws.Range("G11").Value = 2.16
$ws.Range("N11").Value = 2.76
$ws.Range("P11").Value = 1.58
$ws.Range("Q11").Value = 2.64
$ws.Range("T11").Value = 2.18
$ws.Range("U11").Value = 1.79
$ws.Range("W11").Value = 1.86
$ws.Range("X11").Value = 8.6
$ws.Range("AC11").Value = 7
$ws.Range("AD11").Value = 18.5
$ws.Range("AG11").Value = 11
$ws.Range("AM11").Value = 190
$ws.Range("F12").Value = 2.18
$ws.Range("G12").Value = 2.28
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 3.5
$ws.Range("J12").Value = 3.8
$ws.Range("K12").Value = 4.4
$ws.Range("L12").Value = 1.29
$ws.Range("N12").Value = 5
$ws.Range("Q12").Value = 1.6
$ws.Range("T12").Value = 1.57
$ws.Range("U12").Value = 2.46
$ws.Range("V12").Value = 1.41
$ws.Range("W12").Value = 1.78
$ws.Range("X12").Value = 28
$ws.Range("Y12").Value = 22
$ws.Range("AA12").Value = 60
$ws.Range("AC12").Value = 10
$ws.Range("AD12").Value = 18
$ws.Range("AE12").Value = 980
$ws.Range("AF12").Value = 17
$ws.Range("AG12").Value = 12
$ws.Range("AH12").Value = 16.5
$ws.Range("AL12").Value = 980
$ws.Range("AM12").Value = 65
$ws.Range("AN12").Value = 11.5
$ws.Range("AO12").Value = 26
$ws.Range("F13").Value = 2.12
$ws.Range("G13").Value = 2.32
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = 3.45
$ws.Range("J13").Value = 3.85
$ws.Range("K13").Value = 4.4
$ws.Range("L13").Value = 1.2
$ws.Range("O13").Value = 1.12
$ws.Range("S13").Value = 1.94
$ws.Range("T13").Value = 1.46
$ws.Range("V13").Value = 1.41
$ws.Range("W13").Value = 1.74
$ws.Range("AC13").Value = 1000
$ws.Range("AG13").Value = 1000
$ws.Range("AH13").Value = 1000
$ws.Range("AK13").Value = 980
$ws.Range("AN13").Value = 8.6
$ws.Range("AO13").Value = 1000
$ws.Range("H14").Value = 3.05
$ws.Range("N14").Value = 7
$ws.Range("P14").Value = 3.05
$ws.Range("R14").Value = 1.86
$ws.Range("T14").Value = 1.46
$ws.Range("U14").Value = 2.96
$ws.Range("AN14").Value = 9.800000000000001
$ws.Range("F15").Value = 1.87
$ws.Range("G15").Value = 1.91
$ws.Range("I15").Value = 4.3
$ws.Range("N15").Value = 6.6
$ws.Range("S15").Value = 2.08
$ws.Range("T15").Value = 1.5
$ws.Range("V15").Value = 1.31
$ws.Range("W15").Value = 2.08
$ws.Range("Y15").Value = 27
$ws.Range("AD15").Value = 18.5
$ws.Range("AM15").Value = 55
$ws.Range("J16").Value = 3.2
$ws.Range("U16").Value = 2.04
$ws.Range("V16").Value = 1.28
$ws.Range("H17").Value = 2.68
$ws.Range("I17").Value = 2.8
$ws.Range("J17").Value = 3.1
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 2.9
$ws.Range("O17").Value = 1.49
$ws.Range("P17").Value = 1.65
$ws.Range("Q17").Value = 2.48
$ws.Range("R17").Value = 1.22
$ws.Range("S17").Value = 4.9
$ws.Range("T17").Value = 2
$ws.Range("U17").Value = 1.9
$ws.Range("V17").Value = 1.55
$ws.Range("X17").Value = 11.5
$ws.Range("Y17").Value = 9
$ws.Range("Z17").Value = 17.5
$ws.Range("AA17").Value = 48
$ws.Range("AB17").Value = 11.5
$ws.Range("AD17").Value = 13
$ws.Range("AG17").Value = 15
$ws.Range("AH17").Value = 24
$ws.Range("AO17").Value = 44
$ws.Range("H18").Value = 15
$ws.Range("I18").Value = 16.5
$ws.Range("P18").Value = 2.34
$ws.Range("R18").Value = 1.52
$ws.Range("S18").Value = 2.8
$ws.Range("U18").Value = 1.68
$ws.Range("V18").Value = 1.06
$ws.Range("W18").Value = 4.5
$ws.Range("X18").Value = 23
$ws.Range("AA18").Value = 800
$ws.Range("AB18").Value = 8.4
$ws.Range("AC18").Value = 15.5
$ws.Range("AE18").Value = 310
$ws.Range("AI18").Value = 240
$ws.Range("AJ18").Value = 8.800000000000001
$ws.Range("AM18").Value = 260
$ws.Range("AN18").Value = 4.8
$ws.Range("AO18").Value = 410
$ws.Range("F19").Value = 2.84
$ws.Range("G19").Value = 2.9
$ws.Range("H19").Value = 2.78
$ws.Range("I19").Value = 2.82
$ws.Range("J19").Value = 3.35
$ws.Range("N19").Value = 3.35
$ws.Range("O19").Value = 1.39
$ws.Range("P19").Value = 1.79
$ws.Range("Q19").Value = 2.22
$ws.Range("U19").Value = 2.06
$ws.Range("AI19").Value = 50
$ws.Range("AN19").Value = 34
$ws.Range("F18").Value = 1.26
$ws.Range("G18").Value = 1.27
$ws.Range("J18").Value = 6.6
$ws.Range("K18").Value = 6.8
